$d = $word.ActiveDocument

function Insert-SplitRun($ParaRange, $SearchText, $RunsXml) {
    $paraEnd = $ParaRange.End
    $searchRange = $d.Range($ParaRange.Start, $paraEnd - 1)
    $found = $searchRange.Find.Execute($SearchText)
    if (-not $found) {
        throw "Text not found: $SearchText"
    }
    # Re-wrap the found bounds in a *fresh* Range object rather than reusing
    # $searchRange directly: a Range that Find.Execute ran on is left in a
    # state where InsertXML inserts-after instead of replacing its content.
    $target = $d.Range($searchRange.Start, $searchRange.End)
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p>$RunsXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $target.InsertXML($pkg)
}

# --- 1. Title paragraph: split "istema de Tickets soporte técnico para el Centro de Cómputo del ITL" ---
$titleRpr = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='32'/><w:szCs w:val='32'/><w:lang w:val='es-ES'/></w:rPr>"
$titleRuns = "<w:r>$titleRpr<w:t xml:space='preserve'>istema de </w:t></w:r>" +
             "<w:proofErr w:type='gramStart'/>" +
             "<w:r>$titleRpr<w:t>Tickets</w:t></w:r>" +
             "<w:proofErr w:type='gramEnd'/>" +
             "<w:r>$titleRpr<w:t xml:space='preserve'> soporte técnico para el Centro de Cómputo del ITL</w:t></w:r>"
Insert-SplitRun $d.Paragraphs(2).Range "istema de Tickets soporte técnico para el Centro de Cómputo del ITL" $titleRuns

# --- 3. Paragraph "Cuando un usuario termine de llenar su ticket, ..." split around "ticket" ---
$bodyRpr = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>"
$runs9 = "<w:r>$bodyRpr<w:t xml:space='preserve'>Cuando un usuario termine de llenar su </w:t></w:r>" +
         "<w:proofErr w:type='gramStart'/>" +
         "<w:r>$bodyRpr<w:t>ticket</w:t></w:r>" +
         "<w:proofErr w:type='gramEnd'/>" +
         "<w:r>$bodyRpr<w:t xml:space='preserve'>, se le presente una extensión del formulario donde se cuestione sobre la calidad del servicio.</w:t></w:r>"
Insert-SplitRun $d.Paragraphs(9).Range "Cuando un usuario termine de llenar su ticket, se le presente una extensión del formulario donde se cuestione sobre la calidad del servicio." $runs9

# --- 4. Paragraph "Esa información se podría guardar ... Jefe de departamento ..." split around "Jefe" ---
$runs10 = "<w:r>$bodyRpr<w:t xml:space='preserve'>Esa información se podría guardar dentro de la misma API, y se podría gestionar por el </w:t></w:r>" +
          "<w:proofErr w:type='gramStart'/>" +
          "<w:r>$bodyRpr<w:t>Jefe</w:t></w:r>" +
          "<w:proofErr w:type='gramEnd'/>" +
          "<w:r>$bodyRpr<w:t xml:space='preserve'> de departamento, dentro de su apartado de estadísticas, para así poder evaluar que tan satisfactorio o no es el servicio de tickets web.</w:t></w:r>"
Insert-SplitRun $d.Paragraphs(10).Range "Esa información se podría guardar dentro de la misma API, y se podría gestionar por el Jefe de departamento, dentro de su apartado de estadísticas, para así poder evaluar que tan satisfactorio o no es el servicio de tickets web." $runs10

# --- 5. Paragraph "Otra función que podría adherirse al sistema ... ticket de falso a verdadero ..." split around "ticket" ---
$runs11 = "<w:r>$bodyRpr<w:t xml:space='preserve'>Otra función que podría adherirse al sistema es que, el operador al cambiar el estatus del </w:t></w:r>" +
          "<w:proofErr w:type='gramStart'/>" +
          "<w:r>$bodyRpr<w:t>ticket</w:t></w:r>" +
          "<w:proofErr w:type='gramEnd'/>" +
          "<w:r>$bodyRpr<w:t xml:space='preserve'> de falso a verdadero, automáticamente se le mande un correo electrónico al usuario que envió su ticket con la información de que su solicitud ha sido atendida por un operador.</w:t></w:r>"
Insert-SplitRun $d.Paragraphs(11).Range "Otra función que podría adherirse al sistema es que, el operador al cambiar el estatus del ticket de falso a verdadero, automáticamente se le mande un correo electrónico al usuario que envió su ticket con la información de que su solicitud ha sido atendida por un operador." $runs11

# --- 6. Paragraph "Otra función que podría incluirse ..." split around "más" and "ordenamiento" (no proofErr) ---
$runs13 = "<w:r>$bodyRpr<w:t xml:space='preserve'>Otra función que podría incluirse es la de que, en las tablas del sistema, solo se muestren las peticiones que no se han resuelto, ordenar por importancia, u ordenar por </w:t></w:r>" +
          "<w:r>$bodyRpr<w:t>más</w:t></w:r>" +
          "<w:r>$bodyRpr<w:t xml:space='preserve'> antiguos, en general, un método de </w:t></w:r>" +
          "<w:r>$bodyRpr<w:t>ordenamiento</w:t></w:r>" +
          "<w:r>$bodyRpr<w:t xml:space='preserve'> es algo que sería muy útil.</w:t></w:r>"
Insert-SplitRun $d.Paragraphs(13).Range "Otra función que podría incluirse es la de que, en las tablas del sistema, solo se muestren las peticiones que no se han resuelto, ordenar por importancia, u ordenar por mas antiguos, en general, un método de ordenamiento, es algo que sería muy útil." $runs13

# --- 2. Paragraph 7 ("Una función que falto esquematizar..."): add bookmarkStart + bookmarkEnd range ---
# Done last, after all paragraph text has been split into final runs, so the
# bookmark boundaries land on the final run structure (inserting it earlier
# confuses later InsertXML replacements on paragraph 13, which ends with the
# bookmark boundary).
$p7 = $d.Paragraphs(7)
$p13 = $d.Paragraphs(13)
$bookmarkRange = $d.Range($p7.Range.Start, $p13.Range.End - 1)
$d.Bookmarks.Add("_Hlk167752490", $bookmarkRange)

Write-Output "done"
